# ARTIGO4_CriarHistóriasDeUsuários.docx
#
# Commit: "Histórias de Usuários - pequena redefinição de uma história"
#
# The three title/heading lines in the body were each originally split
# across two or more runs that share identical run formatting; the edit
# collapses each of those spans into a single run (and, for the FDD
# line, also removes the now-pointless <w:proofErr/> spell-check
# bookmarks that wrapped the English loan words). A plain Find/Replace
# over the merged text reproduces that run-merge exactly, because Word's
# replace-in-place rewrites the whole matched span as one run.

$d = $word.ActiveDocument

# 1) "ARTIGO TECH #" + "4"  ->  one run "ARTIGO TECH #4"
$d.Content.Find.Execute(
    "ARTIGO TECH #4", $true, $false, $false, $false, $false,
    $true, 1, $false, "ARTIGO TECH #4", 2) | Out-Null

# 2) "Definição de " + "Histórias de Usuários"  ->  one run
$d.Content.Find.Execute(
    "Definição de Histórias de Usuários", $true, $false, $false, $false, $false,
    $true, 1, $false, "Definição de Histórias de Usuários", 2) | Out-Null

# 3) "Técnica FDD – " + "Feature" + " " + "Driven" + " " + "Development" +
#    " (Desenvolvimento Orientado a Funcionalidades)"  ->  one run
#    (also drops the spellcheck proofErr markers around the English words)
$fdd = "Técnica FDD – Feature Driven Development (Desenvolvimento Orientado a Funcionalidades)"
$d.Content.Find.Execute(
    $fdd, $true, $false, $false, $false, $false,
    $true, 1, $false, $fdd, 2) | Out-Null
